$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Station name relabeling (shared strings content) ---
# Red line stations
$ws.Range("B2").Value  = "STATION; Shadyside"
$ws.Range("B3").Value  = "STATION; Herron_Ave"
$ws.Range("B4").Value  = "STATION; Swissville"
$ws.Range("B5").Value  = "STATION; Penn_Station UNDERGROUND"
$ws.Range("B6").Value  = "STATION; Steel_Plaza UNDERGROUND"
$ws.Range("B7").Value  = "STATION; First_Ave; UNDERGROUND"
$ws.Range("B8").Value  = "STATION; Station_Square"
$ws.Range("B9").Value  = "STATION; South_Hills_Junction"

# Green line stations
$ws.Range("B10").Value = "STATION; Pioneer"
$ws.Range("B11").Value = "STATION; Edgebrook"
$ws.Range("B12").Value = "STATION; Station"
$ws.Range("B13").Value = "STATION; Whited"
$ws.Range("B14").Value = "STATION; South_Bank"
$ws.Range("B15").Value = "STATION; Central; UNDERGROUND"
$ws.Range("B16").Value = "STATION; Inglewood; UNDERGROUND"
$ws.Range("B17").Value = "STATION; Glenbury"
$ws.Range("B18").Value = "STATION; Dormont"
$ws.Range("B19").Value = "STATION; Mt_Lebanon"
$ws.Range("B20").Value = "STATION; Poplar"
$ws.Range("B21").Value = "STATION; Castle_Shannon"
$ws.Range("B22").Value = "STATION; Dormont"
$ws.Range("B23").Value = "STATION; Glenbury"
$ws.Range("B24").Value = "STATION; Overbrook; UNDERGROUND"
$ws.Range("B25").Value = "STATION; Inglewood; UNDERGROUND"
$ws.Range("B26").Value = "STATION; Central; UNDERGROUND"

# --- Theme rename ---
try {
  $wb.Theme.Name = "Office Theme 2013 - 2022"
} catch {
}

# --- View / selection update: scroll back to A1 (no topLeftCell override) and select G12 ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("G12").Select() | Out-Null
